$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 5
    6  = 5
    7  = 3
    8  = 5
    9  = 7
    10 = 3
    11 = 2
    12 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
